$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder schema name with a generic template value
$ws.Range("B2").Value = "[DL-MAA20XX-YY]"
$ws.Range("B3").Value = "[DL-MAA20XX-YY]"

# Select B3 to match the resulting cursor position
$ws.Range("B3").Select()
